{"js": "// The commit removes now-stale <w:proofErr/> spell/grammar-check markers\n// (the JSON intro paragraphs and the \"favoriteNumber\"/\"isProgrammer\"/\n// \"hobbies\"/\"friends\"/\"car\" lines each had their text split across several\n// runs only so Word could bracket a misspelling/grammar flag; the text\n// itself is unchanged) and appends a large new \"What is an API?\" section\n// after the JSON example, replacing the trailing tab-only paragraph.\n//\n// Doing this run-by-run would require re-creating dozens of runs/paragraphs\n// by hand and would still have to special-case the <w:tab/> elements (a\n// plain insertText(\"\\t\") turns them into a literal tab character inside\n// <w:t>, not a <w:tab/> run child). Instead, replace the whole body in one\n// shot with the exact target OOXML via insertOoxml, which is the\n// Office.js-native way to splice in a fully-formed run/paragraph tree.\nconst ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>JSON stands for Javascript object notation and is a lightweight data interchange format that replaced XML in the early 2000s.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">Its extremely easy to read/write and is lightweight. </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">It is simply a fomat to transfer data from client to server. </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">JSON can contain Objects, arrays, Booleans, strings, numbers and null. </w:t></w:r></w:p><w:p><w:r><w:t>{</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cname\u201d:\u201dSilas\u201d,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cfavoriteNumber\u201d:10,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cisProgrammer\u201d:true</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201chobbies\u201d:[\u201cWeight Lifting\u201d,\u201dBowling\u201d],</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cfriends\u201d:null,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201ccar\u201d:{</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cmanufacturer\u201d:\u201dFord\u201d,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cmake\u201d:\u201dMustang\u201d</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>What is an API?</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">An API stands for Application Programming Interface, </w:t></w:r><w:r><w:t>which is a set of definitions and protocols for building and integrating application software.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Basically an API is how 2 applications talk to each other. </w:t></w:r></w:p><w:p><w:r><w:t>A type of API that is very popular is known as a RESTful API conforms to the guidelines set by computer scientist roy fielding:</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">1. </w:t></w:r><w:r><w:t>A client-server architecture made up of clients, servers, and resources, with requests managed through HTTP.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">2. </w:t></w:r><w:r><w:t>Stateless client-server communication, meaning no client information is stored between get requests and each request is separate and unconnected.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">3. </w:t></w:r><w:r><w:t>Cacheable data that streamlines client-server interactions.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">3. </w:t></w:r><w:r><w:t>A uniform interface between components so that information is transferred in a standard form. This requires that:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">    a. </w:t></w:r><w:r><w:t>resources requested are identifiable and separate from the representations sent to the client.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">    b. </w:t></w:r><w:r><w:t>resources can be manipulated by the client via the representation they receive because the representation contains enough information to do so.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">    c. </w:t></w:r><w:r><w:t>self-descriptive messages returned to the client have enough information to describe how the client should process it.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">   d. </w:t></w:r><w:r><w:t>hypertext/hypermedia is available, meaning that after accessing a resource the client should be able to use hyperlinks to find all other currently available actions they can take.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">4. </w:t></w:r><w:r><w:t>A layered system that organizes each type of server (those responsible for security, load-balancing, etc.) involved the retrieval of requested information into hierarchies, invisible to the client.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">5. </w:t></w:r><w:r><w:t>Code-on-demand (optional): the ability to send executable code from the server to the client when requested, extending client functionality.</w:t></w:r></w:p><w:p><w:r><w:t>Take a look at the API.</w:t></w:r></w:p><w:p><w:r><w:t>http status codes:</w:t></w:r></w:p><w:p><w:r><w:t>https://developer.mozilla.org/en-US/docs/Web/HTTP/Status</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nconst body = context.document.body;\nbody.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the whole body content (except the final section properties) in one\n# shot via Range.InsertXML so Word's own OOXML writer regenerates the\n# paragraphs/runs (this naturally drops now-stale <w:proofErr/> spell/grammar\n# markers and merges runs whose text did not actually change).\n$newBodyXml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>JSON stands for Javascript object notation and is a lightweight data interchange format that replaced XML in the early 2000s.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">Its extremely easy to read/write and is lightweight. </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">It is simply a fomat to transfer data from client to server. </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">JSON can contain Objects, arrays, Booleans, strings, numbers and null. </w:t></w:r></w:p><w:p><w:r><w:t>{</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cname\u201d:\u201dSilas\u201d,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cfavoriteNumber\u201d:10,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cisProgrammer\u201d:true</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201chobbies\u201d:[\u201cWeight Lifting\u201d,\u201dBowling\u201d],</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cfriends\u201d:null,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201ccar\u201d:{</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cmanufacturer\u201d:\u201dFord\u201d,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>\u201cmake\u201d:\u201dMustang\u201d</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>What is an API?</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">An API stands for Application Programming Interface, </w:t></w:r><w:r><w:t>which is a set of definitions and protocols for building and integrating application software.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Basically an API is how 2 applications talk to each other. </w:t></w:r></w:p><w:p><w:r><w:t>A type of API that is very popular is known as a RESTful API conforms to the guidelines set by computer scientist roy fielding:</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">1. </w:t></w:r><w:r><w:t>A client-server architecture made up of clients, servers, and resources, with requests managed through HTTP.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">2. </w:t></w:r><w:r><w:t>Stateless client-server communication, meaning no client information is stored between get requests and each request is separate and unconnected.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">3. </w:t></w:r><w:r><w:t>Cacheable data that streamlines client-server interactions.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">3. </w:t></w:r><w:r><w:t>A uniform interface between components so that information is transferred in a standard form. This requires that:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">    a. </w:t></w:r><w:r><w:t>resources requested are identifiable and separate from the representations sent to the client.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">    b. </w:t></w:r><w:r><w:t>resources can be manipulated by the client via the representation they receive because the representation contains enough information to do so.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">    c. </w:t></w:r><w:r><w:t>self-descriptive messages returned to the client have enough information to describe how the client should process it.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">   d. </w:t></w:r><w:r><w:t>hypertext/hypermedia is available, meaning that after accessing a resource the client should be able to use hyperlinks to find all other currently available actions they can take.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">4. </w:t></w:r><w:r><w:t>A layered system that organizes each type of server (those responsible for security, load-balancing, etc.) involved the retrieval of requested information into hierarchies, invisible to the client.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">5. </w:t></w:r><w:r><w:t>Code-on-demand (optional): the ability to send executable code from the server to the client when requested, extending client functionality.</w:t></w:r></w:p><w:p><w:r><w:t>Take a look at the API.</w:t></w:r></w:p><w:p><w:r><w:t>http status codes:</w:t></w:r></w:p><w:p><w:r><w:t>https://developer.mozilla.org/en-US/docs/Web/HTTP/Status</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n$d.Content.InsertXML($newBodyXml)\n"}
